$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column F (reuse header formatting from E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Timestamps captured per-row for the new time_taken column
$timestamps = @{
    2  = "2021-10-05 10:52:03.171598"
    3  = "2021-10-05 10:52:03.171608"
    4  = "2021-10-05 10:52:03.171612"
    5  = "2021-10-05 10:52:03.171614"
    6  = "2021-10-05 10:52:03.171617"
    7  = "2021-10-05 10:52:03.171620"
    8  = "2021-10-05 10:52:03.171623"
    9  = "2021-10-05 10:52:03.171625"
    10 = "2021-10-05 10:52:03.171628"
    11 = "2021-10-05 10:52:03.171631"
    12 = "2021-10-05 10:52:03.171633"
    13 = "2021-10-05 10:52:03.171636"
    14 = "2021-10-05 10:52:03.171638"
    15 = "2021-10-05 10:52:03.171641"
    16 = "2021-10-05 10:52:03.171644"
    17 = "2021-10-05 10:52:03.171646"
    18 = "2021-10-05 10:52:03.171649"
    19 = "2021-10-05 10:52:03.171651"
    20 = "2021-10-05 10:52:03.171654"
    21 = "2021-10-05 10:52:03.171656"
    22 = "2021-10-05 10:52:03.171659"
    23 = "2021-10-05 10:52:03.171661"
    24 = "2021-10-05 10:52:03.171664"
    25 = "2021-10-05 10:52:03.171666"
    26 = "2021-10-05 10:52:03.171669"
    27 = "2021-10-05 10:52:03.171672"
    28 = "2021-10-05 10:52:03.171674"
    29 = "2021-10-05 10:52:03.171677"
    30 = "2021-10-05 10:52:03.171679"
}

foreach ($row in 2..30) {
    $ws.Cells.Item($row, 6).Value = $timestamps[$row]
}
